# Refresh market-price-derived profit metrics on each job sheet.
# Values come from an external pricing feed; only data cells are touched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 35042.387
$ws.Range("I64").Value = 74001.36
$ws.Range("J64").Value = 2958.5293
$ws.Range("K64").Value = 74001.36
$ws.Range("L64").Value = 2958.5293
$ws.Range("M64").Value = -73753.36
$ws.Range("N64").Value = -3454.5293
$ws.Range("H67").Value = 35042.387
$ws.Range("I67").Value = 74001.36
$ws.Range("J67").Value = 2958.5293
$ws.Range("K67").Value = 74001.36
$ws.Range("L67").Value = 2958.5293
$ws.Range("M67").Value = -73143.36
$ws.Range("N67").Value = -4674.5293
$ws.Range("H76").Value = 5000.6
$ws.Range("I76").Value = 5000.75
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 5000.75
$ws.Range("L76").Value = 5000
$ws.Range("M76").Value = -4685.75
$ws.Range("N76").Value = -5630
$ws.Range("H79").Value = 5000.6
$ws.Range("I79").Value = 5000.75
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 5000.75
$ws.Range("L79").Value = 5000
$ws.Range("M79").Value = -3908.75
$ws.Range("N79").Value = -7184
$ws.Range("H95").Value = 38000
$ws.Range("J95").Value = 38000
$ws.Range("L95").Value = 38000
$ws.Range("N95").Value = -43492
$ws.Range("H137").Value = 1748.5428
$ws.Range("I137").Value = 1520.3914
$ws.Range("J137").Value = 2185.8333
$ws.Range("K137").Value = 4561.174199999999
$ws.Range("L137").Value = 6557.499899999999
$ws.Range("M137").Value = -2011.174199999999
$ws.Range("N137").Value = -11657.4999
$ws.Range("H138").Value = 2628.236
$ws.Range("I138").Value = 1312.0312
$ws.Range("J138").Value = 3681.2
$ws.Range("K138").Value = 3936.0936
$ws.Range("L138").Value = 11043.6
$ws.Range("M138").Value = 1203.9064
$ws.Range("N138").Value = -21323.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 2416.6667
$ws.Range("J10").Value = 2416.6667
$ws.Range("L10").Value = 2416.6667
$ws.Range("N10").Value = -2756.6667
$ws.Range("H14").Value = 1848.3334
$ws.Range("I14").Value = 445
$ws.Range("J14").Value = 2550
$ws.Range("K14").Value = 445
$ws.Range("L14").Value = 2550
$ws.Range("M14").Value = -270
$ws.Range("N14").Value = -2900
$ws.Range("H46").Value = 4005.5
$ws.Range("J46").Value = 3674
$ws.Range("L46").Value = 3674
$ws.Range("N46").Value = -4312
$ws.Range("H61").Value = 1759.65
$ws.Range("I61").Value = 1379.0625
$ws.Range("K61").Value = 1379.0625
$ws.Range("M61").Value = -1167.0625
$ws.Range("H63").Value = 2200
$ws.Range("J63").Value = 2200
$ws.Range("L63").Value = 2200
$ws.Range("N63").Value = -3572
$ws.Range("H66").Value = 2200
$ws.Range("J66").Value = 2200
$ws.Range("L66").Value = 11000
$ws.Range("N66").Value = -17864
$ws.Range("H101").Value = 49000
$ws.Range("J101").Value = 49000
$ws.Range("L101").Value = 49000
$ws.Range("N101").Value = -55490
$ws.Range("H136").Value = 1759.65
$ws.Range("I136").Value = 1379.0625
$ws.Range("K136").Value = 4137.1875
$ws.Range("M136").Value = -1587.1875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 70768.875
$ws.Range("I86").Value = 93650.164
$ws.Range("J86").Value = 2125
$ws.Range("K86").Value = 93650.164
$ws.Range("L86").Value = 2125
$ws.Range("M86").Value = -92527.164
$ws.Range("N86").Value = -4371
$ws.Range("H89").Value = 70768.875
$ws.Range("I89").Value = 93650.164
$ws.Range("J89").Value = 2125
$ws.Range("K89").Value = 468250.82
$ws.Range("L89").Value = 10625
$ws.Range("M89").Value = -462634.82
$ws.Range("N89").Value = -21857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 8255
$ws.Range("I14").Value = 5010
$ws.Range("J14").Value = 11500
$ws.Range("K14").Value = 5010
$ws.Range("L14").Value = 11500
$ws.Range("M14").Value = -4840
$ws.Range("N14").Value = -11840
$ws.Range("H22").Value = 312.64285
$ws.Range("I22").Value = 234.27272
$ws.Range("K22").Value = 234.27272
$ws.Range("M22").Value = 115.72728
$ws.Range("H31").Value = 38187.125
$ws.Range("I31").Value = 806.4167
$ws.Range("K31").Value = 806.4167
$ws.Range("M31").Value = -511.4167
$ws.Range("H34").Value = 38187.125
$ws.Range("I34").Value = 806.4167
$ws.Range("K34").Value = 806.4167
$ws.Range("M34").Value = -604.4167
$ws.Range("H62").Value = 2357.1428
$ws.Range("J62").Value = 2555.5557
$ws.Range("L62").Value = 2555.5557
$ws.Range("N62").Value = -3803.5557
$ws.Range("H65").Value = 2357.1428
$ws.Range("J65").Value = 2555.5557
$ws.Range("L65").Value = 12777.7785
$ws.Range("N65").Value = -19017.7785
$ws.Range("H68").Value = 17188.412
$ws.Range("J68").Value = 17188.412
$ws.Range("L68").Value = 17188.412
$ws.Range("N68").Value = -18686.412
$ws.Range("H71").Value = 17188.412
$ws.Range("J71").Value = 17188.412
$ws.Range("L71").Value = 51565.236
$ws.Range("N71").Value = -59053.236

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 824546.25
$ws.Range("I131").Value = 521.5833
$ws.Range("J131").Value = 1022312.1
$ws.Range("K131").Value = 1564.7499
$ws.Range("L131").Value = 3066936.3
$ws.Range("M131").Value = 3475.2501
$ws.Range("N131").Value = -3077016.3
$ws.Range("H136").Value = 1431.3
$ws.Range("I136").Value = 1431.3
$ws.Range("K136").Value = 4293.9
$ws.Range("M136").Value = 806.1000000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H51").Value = 40249.5
$ws.Range("J51").Value = 40249.5
$ws.Range("L51").Value = 40249.5
$ws.Range("N51").Value = -41267.5
$ws.Range("H52").Value = 8801.666999999999
$ws.Range("I52").Value = 4030
$ws.Range("J52").Value = 9756
$ws.Range("K52").Value = 4030
$ws.Range("L52").Value = 9756
$ws.Range("M52").Value = -3771
$ws.Range("N52").Value = -10274
$ws.Range("H53").Value = 22331.334
$ws.Range("I53").Value = 7000
$ws.Range("J53").Value = 29997
$ws.Range("K53").Value = 7000
$ws.Range("L53").Value = 29997
$ws.Range("M53").Value = -6369
$ws.Range("N53").Value = -31259
$ws.Range("H70").Value = 64562.47
$ws.Range("I70").Value = 110636.734
$ws.Range("J70").Value = 6201.7334
$ws.Range("K70").Value = 110636.734
$ws.Range("L70").Value = 6201.7334
$ws.Range("M70").Value = -110366.734
$ws.Range("N70").Value = -6741.7334
$ws.Range("H73").Value = 64562.47
$ws.Range("I73").Value = 110636.734
$ws.Range("J73").Value = 6201.7334
$ws.Range("K73").Value = 110636.734
$ws.Range("L73").Value = 6201.7334
$ws.Range("M73").Value = -109700.734
$ws.Range("N73").Value = -8073.7334
$ws.Range("H74").Value = 5531
$ws.Range("J74").Value = 5531
$ws.Range("L74").Value = 5531
$ws.Range("N74").Value = -7403
$ws.Range("H77").Value = 5531
$ws.Range("J77").Value = 5531
$ws.Range("L77").Value = 16593
$ws.Range("N77").Value = -25953
$ws.Range("H80").Value = 167035170
$ws.Range("I80").Value = 200441600
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 200441600
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -200440602
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 167035170
$ws.Range("I83").Value = 200441600
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 1002208000
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -1002203008
$ws.Range("N83").Value = -24984
$ws.Range("H113").Value = 1703
$ws.Range("I113").Value = 1205.5
$ws.Range("J113").Value = 1868.8334
$ws.Range("K113").Value = 1205.5
$ws.Range("L113").Value = 1868.8334
$ws.Range("M113").Value = 964.5
$ws.Range("N113").Value = -6208.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H57").Value = 9833.333000000001
$ws.Range("I57").Value = 500
$ws.Range("J57").Value = 14500
$ws.Range("K57").Value = 500
$ws.Range("L57").Value = 14500
$ws.Range("M57").Value = 66
$ws.Range("N57").Value = -15632
$ws.Range("H136").Value = 1986.8667
$ws.Range("I136").Value = 1600.4445
$ws.Range("J136").Value = 2566.5
$ws.Range("K136").Value = 4801.333500000001
$ws.Range("L136").Value = 7699.5
$ws.Range("M136").Value = -2251.333500000001
$ws.Range("N136").Value = -12799.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 8380.857
$ws.Range("J45").Value = 8944.333000000001
$ws.Range("L45").Value = 8944.333000000001
$ws.Range("N45").Value = -9926.333000000001
$ws.Range("H122").Value = 2792.7778
$ws.Range("I122").Value = 2017.6666
$ws.Range("K122").Value = 6052.9998
$ws.Range("M122").Value = -3602.9998
$ws.Range("H132").Value = 6395.0938
$ws.Range("I132").Value = 3652.9524
$ws.Range("K132").Value = 10958.8572
$ws.Range("M132").Value = -8428.8572
